$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("酒店组播")

# Rename channel "陕西生活" -> "陕西银龄" for its three IP rows (A8:A10)
$ws.Range("A8").Value = "陕西银龄"
$ws.Range("A9").Value = "陕西银龄"
$ws.Range("A10").Value = "陕西银龄"

# Rename channel "陕西公共" -> "陕西秦腔" for its three IP rows (A14:A16)
$ws.Range("A14").Value = "陕西秦腔"
$ws.Range("A15").Value = "陕西秦腔"
$ws.Range("A16").Value = "陕西秦腔"
